# Insert a new weekly data row at row 476 (shifting existing rows 476-508 down to 477-509)
# and populate it with the new observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(476).Insert()

$ws.Cells.Item(476, 1).Value = 9
$ws.Cells.Item(476, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(476, 3).Value = "Metropolitana"
$ws.Cells.Item(476, 4).Value = 44931
$ws.Cells.Item(476, 5).Value = 13
$ws.Cells.Item(476, 6).Value = 100112039
$ws.Cells.Item(476, 7).Value = "Ciboulette"
$ws.Cells.Item(476, 8).Value = "Sin especificar"
$ws.Cells.Item(476, 9).Value = "Primera"
$ws.Cells.Item(476, 10).Value = 340
$ws.Cells.Item(476, 11).Value = 1000
$ws.Cells.Item(476, 12).Value = 1000
$ws.Cells.Item(476, 13).Value = 1000
$ws.Cells.Item(476, 14).Value = "`$/docena de atados"
$ws.Cells.Item(476, 15).Value = "Región Metropolitana"
$ws.Cells.Item(476, 16).Value = 333
$ws.Cells.Item(476, 17).Value = 3
$ws.Cells.Item(476, 18).Value = "Hortaliza"
